$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.280.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").Value = "'1.638.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.73%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "'216.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("D6").Value = "'0.523"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.05%  "

$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("E8").Value = "  -0.17%  "

$ws.Range("D9").Value = "'0.0628"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.14%  "

$ws.Range("D10").Value = "'20.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.39%  "

$ws.Range("D11").Value = "'0.0851"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.20%  "

$ws.Range("D12").Value = "'1.635.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.05%  "

$ws.Range("D13").Value = "'4.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.08%  "

$ws.Range("D14").Value = "'0.549"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.32%  "

$ws.Range("D15").Value = "'65.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.21%  "

$ws.Range("D16").Value = "'27.210.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").Value = "'0.0₃0743"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.51%  "

$ws.Range("D18").Value = "'219.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.21%  "

$ws.Range("D20").Value = "'7.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.19%  "

$ws.Range("D21").Value = "'4.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.46%  "

$ws.Range("E22").Value = "  -5.86%  "

$ws.Range("D23").Value = "'9.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.16%  "

$ws.Range("D24").Value = "'147.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.16%  "

$ws.Range("E25").Value = "  -0.09%  "

$ws.Range("D26").Value = "'7.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.88%  "

$ws.Range("D27").Value = "'0.119"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.21%  "

$ws.Range("D28").Value = "'15.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.60%  "

$ws.Range("D29").Value = "'0.0509"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.22%  "

$ws.Range("E30").Value = "  -0.50%  "

$ws.Range("D31").Value = "'3.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.23%  "

$ws.Range("D32").Value = "'3.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.92%  "

$ws.Range("D33").Value = "'1.338.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.23%  "

$ws.Range("D34").Value = "'1.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.27%  "

$ws.Range("E35").Value = "  -0.38%  "

$ws.Range("E36").Value = "  -0.84%  "

$ws.Range("D37").Value = "'0.549"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("D38").Value = "'0.854"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.67%  "

$ws.Range("E39").Value = "  -0.12%  "

$ws.Range("E40").Value = "  +1.63%  "

$ws.Range("D41").Value = "'0.804"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.70%  "

$ws.Range("D42").Value = "'64.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.02%  "

$ws.Range("D43").Value = "'1.773.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.04%  "

$ws.Range("D44").Value = "'5.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.50%  "

$ws.Range("D45").Value = "'91.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.61%  "

$ws.Range("E46").Value = "  +0.96%  "

$ws.Range("E47").Value = "  -0.05%  "

$ws.Range("D48").Value = "'0.811"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +21.47%  "

$ws.Range("D49").Value = "'0.0515"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.24%  "

$ws.Range("D50").Value = "'0.0992"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.83%  "

$ws.Range("D51").Value = "'7.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.31%  "
